# Generate Report for Handback
#
# The handoff files have now been handed back and are in sync with
# en-US, so this updates the localization-status report to reflect
# that: the per-language "Status" is updated, each row's "Latest
# Target File" / "Latest Handback File" columns are populated (same
# file references as the handoff columns for this sample data), and
# the "Latest Handback DateTime" is stamped.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet ----------------------------------------------------
# The Overview sheet mirrors each language's Status for every file, so it
# also needs to reflect the new handback status.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("B2").Value = $newStatus
$ws.Range("B3").Value = $newStatus

$ws.Range("E2").Value = "56e50667-f644-4794-a1eb-447628c65d1d.md"
$ws.Hyperlinks.Add(
    $ws.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/a6673706d9d9617bece3a7f8d41ddbdf6486976f/e2e/56e50667-f644-4794-a1eb-447628c65d1d.md",
    $null,
    $null,
    "56e50667-f644-4794-a1eb-447628c65d1d.md"
) | Out-Null

$ws.Range("F2").Value = "56e50667-f644-4794-a1eb-447628c65d1d.69015440045a5dd54f97992ae37be4cdf9c899d2.zh-cn.xlf"
$ws.Hyperlinks.Add(
    $ws.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b6d61dada16c79bca4d4797acbbde3870583be73/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/56e50667-f644-4794-a1eb-447628c65d1d.69015440045a5dd54f97992ae37be4cdf9c899d2.zh-cn.xlf",
    $null,
    $null,
    "56e50667-f644-4794-a1eb-447628c65d1d.69015440045a5dd54f97992ae37be4cdf9c899d2.zh-cn.xlf"
) | Out-Null

$ws.Range("G2").Value = "2016-03-04 05:57:44"

$ws.Range("E3").Value = "a38ee76d-c2a0-4161-98f4-99e9098f354a.md"
$ws.Hyperlinks.Add(
    $ws.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/a6673706d9d9617bece3a7f8d41ddbdf6486976f/e2e/a38ee76d-c2a0-4161-98f4-99e9098f354a.md",
    $null,
    $null,
    "a38ee76d-c2a0-4161-98f4-99e9098f354a.md"
) | Out-Null

$ws.Range("F3").Value = "a38ee76d-c2a0-4161-98f4-99e9098f354a.0cd68da6d46ab54d74c9441f133409dda5b7dba6.zh-cn.xlf"
$ws.Hyperlinks.Add(
    $ws.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b6d61dada16c79bca4d4797acbbde3870583be73/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a38ee76d-c2a0-4161-98f4-99e9098f354a.0cd68da6d46ab54d74c9441f133409dda5b7dba6.zh-cn.xlf",
    $null,
    $null,
    "a38ee76d-c2a0-4161-98f4-99e9098f354a.0cd68da6d46ab54d74c9441f133409dda5b7dba6.zh-cn.xlf"
) | Out-Null

$ws.Range("G3").Value = "2016-03-04 05:57:44"

# --- de-de sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("B2").Value = $newStatus
$ws.Range("B3").Value = $newStatus

$ws.Range("E2").Value = "56e50667-f644-4794-a1eb-447628c65d1d.md"
$ws.Hyperlinks.Add(
    $ws.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/a6673706d9d9617bece3a7f8d41ddbdf6486976f/e2e/56e50667-f644-4794-a1eb-447628c65d1d.md",
    $null,
    $null,
    "56e50667-f644-4794-a1eb-447628c65d1d.md"
) | Out-Null

$ws.Range("F2").Value = "56e50667-f644-4794-a1eb-447628c65d1d.69015440045a5dd54f97992ae37be4cdf9c899d2.de-de.xlf"
$ws.Hyperlinks.Add(
    $ws.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aa96c78383c3e6cfd4c88ad5540d6fc9f7ace460/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/56e50667-f644-4794-a1eb-447628c65d1d.69015440045a5dd54f97992ae37be4cdf9c899d2.de-de.xlf",
    $null,
    $null,
    "56e50667-f644-4794-a1eb-447628c65d1d.69015440045a5dd54f97992ae37be4cdf9c899d2.de-de.xlf"
) | Out-Null

$ws.Range("G2").Value = "2016-03-04 05:58:13"

$ws.Range("E3").Value = "a38ee76d-c2a0-4161-98f4-99e9098f354a.md"
$ws.Hyperlinks.Add(
    $ws.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/a6673706d9d9617bece3a7f8d41ddbdf6486976f/e2e/a38ee76d-c2a0-4161-98f4-99e9098f354a.md",
    $null,
    $null,
    "a38ee76d-c2a0-4161-98f4-99e9098f354a.md"
) | Out-Null

$ws.Range("F3").Value = "a38ee76d-c2a0-4161-98f4-99e9098f354a.0cd68da6d46ab54d74c9441f133409dda5b7dba6.de-de.xlf"
$ws.Hyperlinks.Add(
    $ws.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aa96c78383c3e6cfd4c88ad5540d6fc9f7ace460/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a38ee76d-c2a0-4161-98f4-99e9098f354a.0cd68da6d46ab54d74c9441f133409dda5b7dba6.de-de.xlf",
    $null,
    $null,
    "a38ee76d-c2a0-4161-98f4-99e9098f354a.0cd68da6d46ab54d74c9441f133409dda5b7dba6.de-de.xlf"
) | Out-Null

$ws.Range("G3").Value = "2016-03-04 05:58:13"
